$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 259, shifting existing rows 259-349 down to 260-350
$ws.Rows.Item(259).Insert()

# Populate the newly inserted row 259 with the new record
$ws.Range("A259").Value = 10
$ws.Range("B259").Value = "Vega Modelo de Temuco"
$ws.Range("C259").Value = "La Araucanía"
$ws.Range("D259").Value = 45027
$ws.Range("E259").Value = 9
$ws.Range("F259").Value = "Fruta"
$ws.Range("G259").Value = 100103
$ws.Range("H259").Value = "Frutos de hueso (carozo)"
$ws.Range("I259").Value = 100103002
$ws.Range("J259").Value = "Ciruela"
$ws.Range("K259").Value = "Angeleno"
$ws.Range("L259").Value = "Primera"
$ws.Range("M259").Value = 55
$ws.Range("N259").Value = 14000
$ws.Range("O259").Value = 15000
$ws.Range("P259").Value = 14545
$ws.Range("Q259").Value = "$/bandeja 18 kilos granel"
$ws.Range("R259").Value = "Región de O'Higgins"
$ws.Range("S259").Value = 808
$ws.Range("T259").Value = 18
